# Usuarios: "Incluido" -> "Excluido" header, re-sort the user names
# alphabetically in column A and drop the (all "X") column B marks.
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Usuarios")

$ws1.Range("B1").Value = "Excluido"

$names = @(
    "CASTILLO ROMERO, LÍA",
    "CIDÓN HOFFMAN, JAIME",
    "CIMAS CAMPOS, NOIVE",
    "CUADRIELLO GALDÓS, ÁNGELA",
    "CUEVAS RODRIGUEZ, SARA",
    "DE SANTIAGO ABASCAL, GABRIELA",
    "DE SÁDABA IGAREDA, CELIA",
    "GUTIERREZ GUTIERREZ, MARTA",
    "MUÑOZ GARCÍA, JAIME",
    "Pérez González, Docente",
    "REVUELTA DIAZ, CRISTINA",
    "RUIZ CORTÉS, CARLOTA",
    "SAL SARRIA, SAÚL"
)

for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $i + 2
    $ws1.Cells.Item($row, 1).Value = $names[$i]
    $ws1.Cells.Item($row, 2).ClearContents()
}

# Recursos becomes the active sheet / tab, with A3:B8 selected.
$ws2 = $wb.Worksheets.Item("Recursos")
$ws2.Activate()
$ws2.Range("A3:B8").Select()
